# Chris Gayle.xlsx — add matchNo column + populate full innings-by-innings
# batting log (10 rows) for Chris Gayle, and rename the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename the worksheet tab: "Sheet1" -> "Chris Gayle"
# ---------------------------------------------------------------------
$ws.Name = "Chris Gayle"

# ---------------------------------------------------------------------
# 2) Header row (13 columns now, "matchNo" inserted as new column A)
# ---------------------------------------------------------------------
$headers = @(
    "matchNo", "teamName", "batterName", "states", "runs", "balls",
    "fours", "sixes", "sr", "opponentTeamName", "venue", "date", "result"
)

# ---------------------------------------------------------------------
# 3) Data rows (one per innings), matchNo first
# ---------------------------------------------------------------------
$data = @(
    @("37th","Punjab Kings","Chris Gayle","lbw b Rashid Khan","14","17","1","0","82.35","Sunrisers Hyderabad","Sharjah","September 25","Punjab Kings won by 5 runs"),
    @("42nd","Punjab Kings","Chris Gayle","c HH Pandya b Pollard","1","4","0","0","25.00","Mumbai Indians","Abu Dhabi","September 28","Mumbai won by 6 wickets (with 6 balls remaining)"),
    @("8th","Punjab Kings","Chris Gayle","c Jadeja b Chahar","10","10","2","0","100.00","Chennai Super Kings","Wankhede","April 16","Super Kings won by 6 wickets (with 26 balls remaining)"),
    @("11th","Punjab Kings","Chris Gayle","c sub (RV Patel) b Woakes","11","9","0","1","122.22","Delhi Capitals","Wankhede","April 18","Capitals won by 6 wickets (with 10 balls remaining)"),
    @("21st","Punjab Kings","Chris Gayle","c †Karthik b Shivam Mavi","0","1","0","0","0.00","Kolkata Knight Riders","Ahmedabad","April 26","KKR won by 5 wickets (with 20 balls remaining)"),
    @("26th","Punjab Kings","Chris Gayle","c †de Villiers b Sams","46","24","6","2","191.66","Royal Challengers Bangalore","Ahmedabad","April 30","Punjab Kings won by 34 runs"),
    @("29th","Punjab Kings","Chris Gayle","b Rabada","13","9","1","1","144.44","Delhi Capitals","Ahmedabad","May 02","Capitals won by 7 wickets (with 14 balls remaining)"),
    @("14th","Punjab Kings","Chris Gayle","lbw b Rashid Khan","15","17","2","0","88.23","Sunrisers Hyderabad","Chennai","April 21","Sunrisers won by 9 wickets (with 8 balls remaining)"),
    @("17th","Punjab Kings","Chris Gayle","","43","35","5","2","122.85","Mumbai Indians","Chennai","April 23","Punjab Kings won by 9 wickets (with 14 balls remaining)"),
    @("4th","Punjab Kings","Chris Gayle","c Stokes b Parag","40","28","4","2","142.85","Rajasthan Royals","Wankhede","April 12","Punjab Kings won by 4 runs")
)

$lastCol = $headers.Length      # 13 -> column M
$lastRow = 1 + $data.Length     # 11

# Every one of these cells must be stored as TEXT (t="str"), even the
# purely-numeric-looking ones ("14", "17", "82.35", ...), so force the
# whole used range to a text number format before writing any values.
$fullRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$fullRange.NumberFormat = "@"

# Write header row
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Write data rows (skip genuinely empty fields so the cell is left blank,
# matching real Excel semantics where assigning "" clears the cell rather
# than storing a zero-length text value).
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        if ($val -ne "") {
            $ws.Cells.Item($r + 2, $c + 1).Value = $val
        }
    }
}
